$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# The speaker (Hanqing Liu) fills in the "[演讲时间]" (speaking date) field
# in B3. It was typed with a leading apostrophe (quote prefix) so Excel
# stores it as literal text `"03-30-2022"` rather than parsing it, and a
# short-date number format is carried on the cell even though it displays
# as text.
$cell = $ws.Range("B3")
$cell.Value = '''"03-30-2022"'
$cell.NumberFormat = "mm-dd-yy"

# Update the window's view state: scrolled right one column and the
# active selection left on B18 after entry.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B18").Select()
